# Apply updated crypto price/volume data per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.710.68"
$ws.Range("E2").Value = "  +0.69%  "
$ws.Range("D3").Value = "1.951.32"
$ws.Range("E3").Value = "  +2.00%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.22"
$ws.Range("E5").Value = "  +1.26%  "
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4831"
$ws.Range("E7").Value = "  -0.49%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2944"
$ws.Range("E8").Value = "  +1.61%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06819"
$ws.Range("E9").Value = "  +1.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "112.30"
$ws.Range("E10").Value = "  +2.77%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.47"
$ws.Range("E11").Value = "  +0.95%  "
$ws.Range("D12").Value = "1.942.77"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.528"
$ws.Range("E13").Value = "  +5.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07644"
$ws.Range("E14").Value = "  +1.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6926"
$ws.Range("E15").Value = "  +3.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "296.40"
$ws.Range("E16").Value = "  +9.23%  "
$ws.Range("D17").Value = "30.774.03"
$ws.Range("E17").Value = "  +0.91%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.33"
$ws.Range("E18").Value = "  +3.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.701"
$ws.Range("E19").Value = "  +3.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007694"
$ws.Range("E20").Value = "  +2.17%  "
$ws.Range("D21").Value = "2.206.88"
$ws.Range("E21").Value = "  +1.71%  "
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.551"
$ws.Range("E24").Value = "  +2.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.790"
$ws.Range("E25").Value = "  +4.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.96"
$ws.Range("E26").Value = "  +2.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.39"
$ws.Range("E27").Value = "  +1.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.182"
$ws.Range("E28").Value = "  +4.13%  "
$ws.Range("E29").Value = "  +4.23%  "
$ws.Range("E30").Value = "  +2.75%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.776"
$ws.Range("E31").Value = "  +18.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.429"
$ws.Range("E32").Value = "  +7.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05094"
$ws.Range("E33").Value = "  +2.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7801"
$ws.Range("E34").Value = "  +7.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.162"
$ws.Range("E35").Value = "  +2.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02072"
$ws.Range("E36").Value = "  +2.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.734"
$ws.Range("E37").Value = "  +0.73%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.704"
$ws.Range("E38").Value = "  +1.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.040"
$ws.Range("E39").Value = "  +1.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "110.74"
$ws.Range("E40").Value = "  -0.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4454"
$ws.Range("E41").Value = "  +0.78%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8728"
$ws.Range("E42").Value = "  +0.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.925"
$ws.Range("E43").Value = "  +1.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "70.78"
$ws.Range("E44").Value = "  +4.62%  "
$ws.Range("E45").Value = "  +0.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.377"
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.497"
$ws.Range("E47").Value = "  +3.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "48.68"
$ws.Range("E48").Value = "  +2.75%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1251"
$ws.Range("E49").Value = "  +0.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "35.49"
$ws.Range("E50").Value = "  +2.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.2546"
$ws.Range("E51").Value = "  +2.98%  "
